$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM Simple")

# --- Row 54: new row, the fuse holder (designator XF1) that used to live in row 11 ---
# (populated first so the new shared-string "XF1" is registered before "Fusible ...")
$ws.Range("A54").Formula = '=ROW(A54) - ROW($A$4)'
$ws.Range("B54").Value = "XF1"
$ws.Range("C54").Value = "FUSE HOLDER 5x20MM THT"
$ws.Range("D54").Value = "PN-533361"
$ws.Range("E54").Value = 1
$ws.Range("F54").Value = "Farnell"
$ws.Range("G54").Value = 3517015
$ws.Range("H54").Value = 1.11
$ws.Range("I54").Formula = '=$C$3*E54'
$ws.Range("J54").Formula = '=H54*I54'
$ws.Range("K54").Value = $true

# --- Row 11: was "FUSE HOLDER 5x20MM THT" (F1), becomes "Fusible FST 2A temporisé" ---
$ws.Range("C11").Value = "Fusible FST 2A temporisé"
$ws.Range("D11").Value = 261.476
$ws.Range("D11").HorizontalAlignment = -4131
$ws.Range("F11").Value = "ETML"
$ws.Range("G11").Value = 261.476
$ws.Range("H11").Value = 0.33

$wb.Save()
